$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.230855
$ws.Range("H2").Value = 0.692565
$ws.Range("I2").Value = 0.06377305075821572
$ws.Range("J2").Value = 0.06377305075821572
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.160945333333332
$ws.Range("N2").Value = 24.482836
$ws.Range("O2").Value = 0.491005088714322
$ws.Range("P2").Value = 0.491005088714322
$ws.Range("Q2").Value = 1.883995034926667
$ws.Range("R2").Value = 16.95595531434
$ws.Range("S2").Value = 0.03131289244512067
$ws.Range("T2").Value = 0.03131289244512067

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.230855
$ws.Range("H3").Value = 0.692565
$ws.Range("I3").Value = 0.06377305075821572
$ws.Range("J3").Value = 0.06377305075821572
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.477462333333333
$ws.Range("N3").Value = 7.432386999999999
$ws.Range("O3").Value = 0.1490570715865667
$ws.Range("P3").Value = 0.1490570715865668
$ws.Range("Q3").Value = 0.5719345669616666
$ws.Range("R3").Value = 5.147411102654999
$ws.Range("S3").Value = 0.009505824192161115
$ws.Range("T3").Value = 0.009505824192161116

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.230855
$ws.Range("H4").Value = 0.692565
$ws.Range("I4").Value = 0.06377305075821572
$ws.Range("J4").Value = 0.06377305075821572
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.982489999999999
$ws.Range("N4").Value = 17.94747
$ws.Range("O4").Value = 0.3599378396991113
$ws.Range("P4").Value = 0.3599378396991114
$ws.Range("Q4").Value = 1.38108772895
$ws.Range("R4").Value = 12.42978956055
$ws.Range("S4").Value = 0.02295433412093394
$ws.Range("T4").Value = 0.02295433412093394

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.900731333333333
$ws.Range("H5").Value = 8.702194
$ws.Range("I5").Value = 0.8013189515350044
$ws.Range("J5").Value = 0.8013189515350045
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.160945333333332
$ws.Range("N5").Value = 24.482836
$ws.Range("O5").Value = 0.491005088714322
$ws.Range("P5").Value = 0.491005088714322
$ws.Range("Q5").Value = 23.67270983802044
$ws.Range("R5").Value = 213.054388542184
$ws.Range("S5").Value = 0.3934516828869123
$ws.Range("T5").Value = 0.3934516828869123

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.900731333333333
$ws.Range("H6").Value = 8.702194
$ws.Range("I6").Value = 0.8013189515350044
$ws.Range("J6").Value = 0.8013189515350045
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.477462333333333
$ws.Range("N6").Value = 7.432386999999999
$ws.Range("O6").Value = 0.1490570715865667
$ws.Range("P6").Value = 0.1490570715865668
$ws.Range("Q6").Value = 7.18645261745311
$ws.Range("R6").Value = 64.678073557078
$ws.Range("S6").Value = 0.1194422563226257
$ws.Range("T6").Value = 0.1194422563226258

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.900731333333333
$ws.Range("H7").Value = 8.702194
$ws.Range("I7").Value = 0.8013189515350044
$ws.Range("J7").Value = 0.8013189515350045
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.982489999999999
$ws.Range("N7").Value = 17.94747
$ws.Range("O7").Value = 0.3599378396991113
$ws.Range("P7").Value = 0.3599378396991114
$ws.Range("Q7").Value = 17.35359619435333
$ws.Range("R7").Value = 156.18236574918
$ws.Range("S7").Value = 0.2884250123254664
$ws.Range("T7").Value = 0.2884250123254665

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efna5"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4883596666666667
$ws.Range("H8").Value = 1.465079
$ws.Range("I8").Value = 0.1349079977067798
$ws.Range("J8").Value = 0.1349079977067798
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.160945333333332
$ws.Range("N8").Value = 24.482836
$ws.Range("O8").Value = 0.491005088714322
$ws.Range("P8").Value = 0.491005088714322
$ws.Range("Q8").Value = 3.985476542671555
$ws.Range("R8").Value = 35.869288884044
$ws.Range("S8").Value = 0.06624051338228897
$ws.Range("T8").Value = 0.06624051338228897

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efna5"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4883596666666667
$ws.Range("H9").Value = 1.465079
$ws.Range("I9").Value = 0.1349079977067798
$ws.Range("J9").Value = 0.1349079977067798
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.477462333333333
$ws.Range("N9").Value = 7.432386999999999
$ws.Range("O9").Value = 0.1490570715865667
$ws.Range("P9").Value = 0.1490570715865668
$ws.Range("Q9").Value = 1.209892679285889
$ws.Range("R9").Value = 10.889034113573
$ws.Range("S9").Value = 0.02010899107177986
$ws.Range("T9").Value = 0.02010899107177986

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efna5"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4883596666666667
$ws.Range("H10").Value = 1.465079
$ws.Range("I10").Value = 0.1349079977067798
$ws.Range("J10").Value = 0.1349079977067798
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.982489999999999
$ws.Range("N10").Value = 17.94747
$ws.Range("O10").Value = 0.3599378396991113
$ws.Range("P10").Value = 0.3599378396991114
$ws.Range("Q10").Value = 2.921606822236666
$ws.Range("R10").Value = 26.29446140013
$ws.Range("S10").Value = 0.04855849325271099
$ws.Range("T10").Value = 0.04855849325271099
